# Update DM integration test fixture
#
# - Bold the header row on each of the 3 sheets (CodeSchemes, Codes, Extensions)
# - Widen the columns to accommodate the bold header text
# - Update the ID (GUID) of the CodeSchemes row to a new value

$wb = $excel.ActiveWorkbook

# --- CodeSchemes sheet --------------------------------------------------
$wsCodeSchemes = $wb.Worksheets.Item("CodeSchemes")
$wsCodeSchemes.Range("A1:N1").Font.Bold = $true

$wsCodeSchemes.Columns.Item(1).ColumnWidth = 34.5
$wsCodeSchemes.Columns.Item(2).ColumnWidth = 17.5
$wsCodeSchemes.Columns.Item(3).ColumnWidth = 25.65
$wsCodeSchemes.Columns.Item(4).ColumnWidth = 22.35
$wsCodeSchemes.Columns.Item(5).ColumnWidth = 14.1
$wsCodeSchemes.Columns.Item(6).ColumnWidth = 19.05
$wsCodeSchemes.Columns.Item(7).ColumnWidth = 21.3
$wsCodeSchemes.Columns.Item(8).ColumnWidth = 19.05
$wsCodeSchemes.Columns.Item(9).ColumnWidth = 20.75
$wsCodeSchemes.Columns.Item(10).ColumnWidth = 24.05
$wsCodeSchemes.Columns.Item(11).ColumnWidth = 19.05
$wsCodeSchemes.Columns.Item(12).ColumnWidth = 15.8
$wsCodeSchemes.Columns.Item(13).ColumnWidth = 20.75
$wsCodeSchemes.Columns.Item(14).ColumnWidth = 27.35

# Update the row-2 ID value (new GUID)
$wsCodeSchemes.Range("A2").Value = "ed5c8cc0-3299-463d-9474-56242187e817"

# --- Codes sheet ----------------------------------------------------------
$wsCodes = $wb.Worksheets.Item("Codes")
$wsCodes.Range("A1:J1").Font.Bold = $true

$wsCodes.Columns.Item(1).ColumnWidth = 5.9
$wsCodes.Columns.Item(2).ColumnWidth = 17.5
$wsCodes.Columns.Item(3).ColumnWidth = 15.8
$wsCodes.Columns.Item(4).ColumnWidth = 14.1
$wsCodes.Columns.Item(5).ColumnWidth = 15.8
$wsCodes.Columns.Item(6).ColumnWidth = 19.05
$wsCodes.Columns.Item(7).ColumnWidth = 20.75
$wsCodes.Columns.Item(8).ColumnWidth = 24.05
$wsCodes.Columns.Item(9).ColumnWidth = 19.05
$wsCodes.Columns.Item(10).ColumnWidth = 15.8

# --- Extensions sheet -------------------------------------------------------
$wsExtensions = $wb.Worksheets.Item("Extensions")
$wsExtensions.Range("A1:I1").Font.Bold = $true

$wsExtensions.Columns.Item(1).ColumnWidth = 5.9
$wsExtensions.Columns.Item(2).ColumnWidth = 17.5
$wsExtensions.Columns.Item(3).ColumnWidth = 14.1
$wsExtensions.Columns.Item(4).ColumnWidth = 24.05
$wsExtensions.Columns.Item(5).ColumnWidth = 15.8
$wsExtensions.Columns.Item(6).ColumnWidth = 19.05
$wsExtensions.Columns.Item(7).ColumnWidth = 19.05
$wsExtensions.Columns.Item(8).ColumnWidth = 15.8
$wsExtensions.Columns.Item(9).ColumnWidth = 24.05
